$d = $word.ActiveDocument

# --- Locate the "Mouse" section bullet that ends with 按键状态：按下、释放 ---
$rng = $d.Content
$found = $rng.Find.Execute("按键状态：按下、释放", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the '按键状态：按下、释放' paragraph"
}
$srcPara = $rng.Paragraphs(1)

# --- Insert a brand-new list paragraph right after it; InsertParagraphAfter ---
# --- on a ListParagraph/numPr paragraph carries the same bullet formatting ---
$srcPara.Range.InsertParagraphAfter()
$newPara = $srcPara.Next()
$newParaRange = $newPara.Range

# Type the new bullet text with a one-character placeholder appended.
# The placeholder lets the _GoBack bookmark be dropped at a real run
# boundary (the runtime can't collapse a bookmark on the very last
# character position right before a paragraph mark), after which the
# placeholder is deleted, leaving the bookmark sitting immediately
# after the real text and before the paragraph mark - matching where
# Word itself leaves the hidden _GoBack bookmark after the last edit.
$newParaRange.InsertBefore("鼠标事件：按下、释放Z")

$freshRange = $newPara.Range
$placeholderPos = $freshRange.End - 2

# Word keeps a single _GoBack bookmark; adding it here automatically
# removes it from its old location (right after 假定文本与以).
$bmRange = $d.Range($placeholderPos, $placeholderPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$placeholderRange = $d.Range($placeholderPos, $placeholderPos + 1)
$placeholderRange.Delete()

Write-Output "Inserted 鼠标事件：按下、释放 bullet and relocated _GoBack bookmark"
